$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''246.94'
$ws.Range("E2").Value = '''0.65%'
$ws.Range("D3").Value = '''26.13'
$ws.Range("E3").Value = '''3.85%'
$ws.Range("D4").Value = '''5.089'
$ws.Range("E4").Value = '''0.87%'
$ws.Range("D5").Value = '''0.05595'
$ws.Range("E5").Value = '''-0.20%'
$ws.Range("E6").Value = '''-1.16%'
$ws.Range("D7").Value = '''0.8140'
$ws.Range("E7").Value = '''0.03%'
$ws.Range("D8").Value = '''0.8438'
$ws.Range("E8").Value = '''0.23%'
$ws.Range("D9").Value = '''0.06995'
$ws.Range("E9").Value = '''0.63%'
$ws.Range("D10").Value = '''0.02818'
$ws.Range("E10").Value = '''-0.55%'
$ws.Range("D11").Value = '''0.09392'
$ws.Range("E11").Value = '''-0.11%'
$ws.Range("D12").Value = '''0.001523'
$ws.Range("E12").Value = '''0.79%'
$ws.Range("B13").Value = 'TigerCash'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D13").Value = '''0.006130'
$ws.Range("E13").Value = '''-1.82%'
$ws.Range("B14").Value = 'LEO'
$ws.Range("C14").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D14").Value = '''3.609'
$ws.Range("E14").Value = '''3.15%'
$ws.Range("B15").Value = 'GateToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D15").Value = '''3.022'
$ws.Range("E15").Value = '''0.08%'
$ws.Range("B16").Value = 'BTSEToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D16").Value = '''2.055'
$ws.Range("E16").Value = '''-1.73%'
$ws.Range("B17").Value = 'BitpandaEcosystemToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D17").Value = '''0.3113'
$ws.Range("E17").Value = '''-2.29%'
$ws.Range("B18").Value = 'WazirX'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D18").Value = '''0.1337'
$ws.Range("E18").Value = '''-0.08%'
$ws.Range("B19").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C19").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D19").Value = '''0.03208'
$ws.Range("E19").Value = '''-1.17%'
$ws.Range("B20").Value = 'ProBitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D20").Value = '''0.1276'
$ws.Range("E20").Value = '''-1.28%'
$ws.Range("B21").Value = 'MCDex'
$ws.Range("C21").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D21").Value = '''3.748'
$ws.Range("E21").Value = '''0.25%'
$ws.Range("B22").Value = 'CoinExToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D22").Value = '''0.04651'
$ws.Range("E22").Value = '''-0.77%'
$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D23").Value = '''0.1376'
$ws.Range("E23").Value = '''0.40%'
$ws.Range("B24").Value = 'One'
$ws.Range("C24").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D24").Value = '''0.0005997'
$ws.Range("E24").Value = '''0.60%'
$ws.Range("E25").Value = '''0.11%'
$ws.Range("D26").Value = '''0.004569'
$ws.Range("E26").Value = '''0.87%'
$ws.Range("D27").Value = '''0.00009600'
$ws.Range("E27").Value = '''-1.05%'
$ws.Range("D28").Value = '''0.0001940'
$ws.Range("E28").Value = '''0.00%'
$ws.Range("D40").Value = '''0.03658'
$ws.Range("E40").Value = '''-0.09%'
$ws.Range("D41").Value = '''0.1362'
$ws.Range("E41").Value = '''29.77%'
$ws.Range("B42").Value = 'KickToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D42").Value = '''0.006150'
$ws.Range("E42").Value = '''-0.66%'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = '''0.002617'
$ws.Range("E43").Value = '''-4.13%'
$ws.Range("D44").Value = '''0.008031'
$ws.Range("E44").Value = '''-1.55%'
$ws.Range("D45").Value = '''0.00005381'
$ws.Range("E45").Value = '''1.57%'
$ws.Range("D46").Value = '''0.00000000751'
$ws.Range("E46").Value = '''0.10%'
$ws.Range("D47").Value = '''0.1452'
$ws.Range("E47").Value = '''-19.36%'
$ws.Range("D48").Value = '''0.002427'
$ws.Range("E48").Value = '''20.39%'
$ws.Range("D49").Value = '''0.00002103'
$ws.Range("E49").Value = '''0.11%'
$ws.Range("D50").Value = '''0.0002003'
$ws.Range("E50").Value = '''0.10%'
